# Fix the misspelled "African America" label to "African American"
# on the "ethnicities" worksheet (header cell C1), and update the
# last active selection to D9, matching the saved file's view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ethnicities")

$ws.Range("C1").Value = "African American"

$ws.Activate()
$ws.Range("D9").Select()
